$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (column D) and volume/1h-change (column E) for each coin row.
# A leading apostrophe forces plain-decimal-looking prices to stay text
# (matching the original inline-string cell type) instead of being
# auto-converted to a binary-float number by Excel.
$ws.Cells.Item(2, 4).Value = "25.963.68"
$ws.Cells.Item(2, 5).Value = "  +0.14%  "
$ws.Cells.Item(3, 4).Value = "1.638.38"
$ws.Cells.Item(3, 5).Value = "  -0.57%  "
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = "  -0.74%  "
$ws.Cells.Item(5, 4).Value = "'214.93"
$ws.Cells.Item(5, 5).Value = "  -0.38%  "
$ws.Cells.Item(6, 4).Value = "'0.5135"
$ws.Cells.Item(6, 5).Value = "  +0.64%  "
$ws.Cells.Item(7, 5).Value = "  -0.46%  "
$ws.Cells.Item(8, 4).Value = "'0.2575"
$ws.Cells.Item(8, 5).Value = "  -0.06%  "
$ws.Cells.Item(9, 4).Value = "'0.06349"
$ws.Cells.Item(9, 5).Value = "  -1.13%  "
$ws.Cells.Item(10, 4).Value = "'19.75"
$ws.Cells.Item(10, 5).Value = "  +0.16%  "
$ws.Cells.Item(11, 4).Value = "'0.07769"
$ws.Cells.Item(11, 5).Value = "  -0.02%  "
$ws.Cells.Item(12, 4).Value = "'4.274"
$ws.Cells.Item(12, 5).Value = "  -0.94%  "
$ws.Cells.Item(13, 4).Value = "1.626.92"
$ws.Cells.Item(13, 5).Value = "  -1.61%  "
$ws.Cells.Item(14, 4).Value = "'0.5458"
$ws.Cells.Item(14, 5).Value = "  -0.44%  "
$ws.Cells.Item(15, 4).Value = "0.0₅7741"
$ws.Cells.Item(15, 5).Value = "  -1.94%  "
$ws.Cells.Item(16, 4).Value = "'64.35"
$ws.Cells.Item(16, 5).Value = "  -0.95%  "
$ws.Cells.Item(17, 4).Value = "25.978.16"
$ws.Cells.Item(17, 5).Value = "  -0.05%  "
$ws.Cells.Item(18, 4).Value = "'1.002"
$ws.Cells.Item(18, 5).Value = "  -0.39%  "
$ws.Cells.Item(19, 4).Value = "'197.04"
$ws.Cells.Item(19, 5).Value = "  -0.46%  "
$ws.Cells.Item(20, 4).Value = "'4.430"
$ws.Cells.Item(20, 5).Value = "  -0.17%  "
$ws.Cells.Item(21, 4).Value = "'9.921"
$ws.Cells.Item(21, 5).Value = "  -1.20%  "
$ws.Cells.Item(22, 4).Value = "'6.079"
$ws.Cells.Item(22, 5).Value = "  +0.21%  "
$ws.Cells.Item(25, 4).Value = "'142.28"
$ws.Cells.Item(25, 5).Value = "  +1.15%  "
$ws.Cells.Item(26, 4).Value = "'0.1232"
$ws.Cells.Item(26, 5).Value = "  +7.25%  "
$ws.Cells.Item(27, 4).Value = "'6.836"
$ws.Cells.Item(27, 5).Value = "  -0.95%  "
$ws.Cells.Item(28, 4).Value = "'15.62"
$ws.Cells.Item(28, 5).Value = "  -0.87%  "
$ws.Cells.Item(29, 4).Value = "'1.238"
$ws.Cells.Item(29, 5).Value = "  -0.21%  "
$ws.Cells.Item(30, 4).Value = "'0.04840"
$ws.Cells.Item(30, 5).Value = "  -3.46%  "
$ws.Cells.Item(31, 4).Value = "'3.271"
$ws.Cells.Item(31, 5).Value = "  -0.35%  "
$ws.Cells.Item(32, 4).Value = "'3.209"
$ws.Cells.Item(32, 5).Value = "  +0.03%  "
$ws.Cells.Item(33, 4).Value = "'1.534"
$ws.Cells.Item(33, 5).Value = "  -0.84%  "
$ws.Cells.Item(34, 5).Value = "  +0.32%  "
$ws.Cells.Item(35, 4).Value = "'0.9121"
$ws.Cells.Item(35, 5).Value = "  +1.75%  "
$ws.Cells.Item(36, 4).Value = "'2.568"
$ws.Cells.Item(36, 5).Value = "  -0.72%  "
$ws.Cells.Item(37, 4).Value = "'0.5540"
$ws.Cells.Item(37, 5).Value = "  +0.07%  "
$ws.Cells.Item(38, 4).Value = "1.101.94"
$ws.Cells.Item(38, 5).Value = "  -2.66%  "
$ws.Cells.Item(39, 5).Value = "  +0.11%  "
$ws.Cells.Item(40, 4).Value = "'1.002"
$ws.Cells.Item(40, 5).Value = "  -0.48%  "
$ws.Cells.Item(41, 4).Value = "'2.523"
$ws.Cells.Item(41, 5).Value = "  -1.47%  "
$ws.Cells.Item(42, 4).Value = "'5.555"
$ws.Cells.Item(42, 5).Value = "  -1.86%  "
$ws.Cells.Item(43, 4).Value = "'0.8048"
$ws.Cells.Item(43, 5).Value = "  -1.37%  "
$ws.Cells.Item(44, 4).Value = "'99.13"
$ws.Cells.Item(44, 5).Value = "  -0.77%  "
$ws.Cells.Item(45, 4).Value = "0.0₈120"
$ws.Cells.Item(45, 5).Value = "  -3.93%  "
$ws.Cells.Item(46, 4).Value = "1.779.73"
$ws.Cells.Item(46, 5).Value = "  -0.25%  "
$ws.Cells.Item(47, 4).Value = "'0.4536"
$ws.Cells.Item(47, 5).Value = "  -0.21%  "
$ws.Cells.Item(48, 4).Value = "'55.04"
$ws.Cells.Item(48, 5).Value = "  -0.58%  "
$ws.Cells.Item(49, 4).Value = "'0.9971"
$ws.Cells.Item(49, 5).Value = "  -0.89%  "
$ws.Cells.Item(50, 4).Value = "'0.05209"
$ws.Cells.Item(50, 5).Value = "  +2.22%  "
$ws.Cells.Item(51, 4).Value = "'7.484"
$ws.Cells.Item(51, 5).Value = "  +0.99%  "

# Rows 23 and 24 swap coins: Toncoin <-> BinanceUSD (with refreshed price/volume)
$ws.Cells.Item(23, 2).Value = "BinanceUSD"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(23, 4).Value = "'1.003"
$ws.Cells.Item(23, 5).Value = "  -0.51%  "

$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(24, 4).Value = "'1.933"
$ws.Cells.Item(24, 5).Value = "  +4.29%  "
